$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates: student changed to LIN, YU-TING; add "Assisted on" note ---
$ws.Range("B3").Value = "LIN, YU-TING"
$ws.Range("F3").Value = "Lab 1 questions"

# --- New row 4: KAYLA S REVELLE, 01/29/2016, 3:30 PM-ish, At 3D47, HW1+Lab2 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "KAYLA S REVELLE"

$ws.Range("C4").Formula = "=""01/29/2016"""
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("D4").Value = 0.64583333333333337
$ws.Range("E4").Value = "At 3D47"
$ws.Range("F4").Value = "Lab 2 & Homework 1 questions"

# --- New row 5: LIN, YU-TING, 01/29/2016, At 3D47, HW1+Lab2 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "LIN, YU-TING"

$ws.Range("C5").Formula = "=""01/29/2016"""
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$ws.Range("D5").Value = 0.6875
$ws.Range("E5").Value = "At 3D47"
$ws.Range("F5").Value = "Lab 2 & Homework 1 questions"

# Apply correct formatting (thin border) on the new rows to match existing data rows,
# then restore the Time-format style on column D (numFmtId 18, as used by D2/D3).
$ws.Range("A4:F5").Borders.LineStyle = 1
$ws.Range("D2:D3").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)

# Un-minimize the workbook window
$wb.Windows.Item(1).WindowState = -4143

# Move the active selection to B3
$ws.Range("B3").Select()
